$d = $word.ActiveDocument

# Locate the run containing "County Jail: PS   EM;" so we can append the
# new "Victim's Attorney" runs right after it, inside the same paragraph.
$found = $d.Content.Find.Execute("County Jail: PS   EM;", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)

$searchRange = $d.Content
$searchRange.Find.Execute("County Jail: PS   EM;", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)

# Collapse the found range to its end point, right after the semicolon.
$insertPoint = $d.Range($searchRange.End, $searchRange.End)

# Insert the separating space as its own run.
$insertPoint.InsertAfter(" ")
$insertPoint.Font.Name = "Palatino Linotype"
$insertPoint.Font.Size = 8

# Move past the space we just inserted, then insert the new sentence as
# its own run with the same formatting.
$afterSpace = $d.Range($insertPoint.End, $insertPoint.End)
$afterSpace.InsertAfter("Victim’s Attorney (if applicable): PS   OS   EM")
$afterSpace.Font.Name = "Palatino Linotype"
$afterSpace.Font.Size = 8
